$d = $word.ActiveDocument

# --- Change 1 ---
# " button. You'll see these " + bookmark(_GoBack) + "documents as well as other files"
# becomes a single run " button. You'll see these documents as well as other files"
# and the _GoBack bookmark is removed.
$old1 = " button. You" + [char]8217 + "ll see these documents as well as other files"
$new1 = " button. You" + [char]8217 + "ll see these documents as well as other files"
$rng1 = $d.Content
$rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Change 2 ---
# Split the sentence after "Update": keep " on the right). " in the paragraph,
# drop the trailing sentence about viewing the original PDF (it moves elsewhere).
$old2 = " on the right). If you viewed the original PDF and then you try viewing it again after your update, you may have to click on the refresh button in the browser, since the browser may be showing you the original PDF which it cached when you looked at it the first time."
$new2 = " on the right). "
$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# Find the now-empty paragraph right after that one and put the _GoBack bookmark there.
$rng3 = $d.Content
$rng3.Find.Execute(" on the right). ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$updatePara = $rng3.Paragraphs(1)
$nextPara = $updatePara.Next()
$d.Bookmarks.Add("_GoBack", $nextPara.Range) | Out-Null

# --- Change 3 ---
# Add two new paragraphs after the last picture (the one following the
# "Update ... on the right)." paragraph): an empty paragraph, then a
# paragraph holding the sentence removed in Change 2.
$picPara = $nextPara.Next()
$afterPicRange = $picPara.Range
$afterPicRange.Collapse(0)
$afterPicRange.InsertParagraphAfter()
$afterPicRange.Collapse(0)
$afterPicRange.InsertParagraphAfter()
$afterPicRange.Collapse(0)
$newTextPara = $picPara.Next().Next()
$newTextPara.Range.InsertBefore("If you viewed the original PDF and then you try viewing it again after your update, you may have to click on the refresh button in the browser, since the browser may be showing you the original PDF which it cached when you looked at it the first time.")

Write-Output "done"
